$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Sheet1" worksheet at the end, with the AIRR-format field list ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet1"

$fields = @("sequence_id","sequence","locus","stop_codon","vj_in_frame","productive","rev_comp","complete_vdj","v_call","d_call","j_call","sequence_alignment","germline_alignment","sequence_alignment_aa","germline_alignment_aa","v_alignment_start","v_alignment_end","d_alignment_start","d_alignment_end","j_alignment_start","j_alignment_end","v_sequence_alignment","v_sequence_alignment_aa","v_germline_alignment","v_germline_alignment_aa","d_sequence_alignment","d_sequence_alignment_aa","d_germline_alignment","d_germline_alignment_aa","j_sequence_alignment","j_sequence_alignment_aa","j_germline_alignment","j_germline_alignment_aa","fwr1","fwr1_aa","cdr1","cdr1_aa","fwr2","fwr2_aa","cdr2","cdr2_aa","fwr3","fwr3_aa","fwr4","fwr4_aa","cdr3","cdr3_aa","junction","junction_length","junction_aa","junction_aa_length","v_score","d_score","j_score","v_cigar","d_cigar","j_cigar","v_support","d_support","j_support","v_identity","d_identity","j_identity","v_sequence_start","v_sequence_end","v_germline_start","v_germline_end","d_sequence_start","d_sequence_end","d_germline_start","d_germline_end","j_sequence_start","j_sequence_end","j_germline_start","j_germline_end","fwr1_start","fwr1_end","cdr1_start","cdr1_end","fwr2_start","fwr2_end","cdr2_start","cdr2_end","fwr3_start","fwr3_end","fwr4_start","fwr4_end","cdr3_start","cdr3_end","np1","np1_length","np2","np2_length")

$n = $fields.Length
$arr = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $i
    $arr[$i,1] = $fields[$i]
}
$newSheet.Range("A1:B" + $n).Value = $arr

$newSheet.Columns.Item(1).ColumnWidth = 17
$newSheet.Columns.Item(2).ColumnWidth = 51

# --- 2. Update the "IgBlast" sheet CDR3-related field mapping rows ---
$ws = $wb.Worksheets.Item("IgBlast")

# CDR3-related rows (Sequence/GermlineSequence numeric col index bumped, CDR3* cols switched
# from the placeholder "function" marker to explicit numeric/col-letter codes)
$ws.Range("F81").Value = 11
$ws.Range("F82").Value = 12
$ws.Range("F83").Value = 45
$ws.Range("F84").Value = 46
$ws.Range("F85").Value = 46
$ws.Range("F86").Value = "87-63"
$ws.Range("F87").Value = "88-63"

# TotalMutations / Mutation list rows: "later" placeholder resolved to "function"
$ws.Range("F59").Value = "function"
$ws.Range("F98").Value = "function"
$ws.Range("F99").Value = "function"

# Keep "IgBlast" as the active/selected tab (as it was before the edit)
$ws.Activate()
$ws.Range("D34").Select()
